$wb = $excel.ActiveWorkbook

# ALC row 15
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 3943.5637
$ws.Range("I15").Value = 3943.5637
$ws.Range("K15").Value = 11830.6911
$ws.Range("M15").Value = -11661.6911

# ALC row 47
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H47").Value = 0
$ws.Range("I47").Value = 0
$ws.Range("K47").Value = 0
$ws.Range("M47").ClearContents()

# ALC row 62
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 22226984
$ws.Range("J62").Value = 2006
$ws.Range("L62").Value = 2006
$ws.Range("N62").Value = -3254

# ALC row 65
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 22226984
$ws.Range("J65").Value = 2006
$ws.Range("L65").Value = 10030
$ws.Range("N65").Value = -16270

# ALC row 95
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H95").Value = 10623
$ws.Range("J95").Value = 10623
$ws.Range("L95").Value = 10623
$ws.Range("N95").Value = -16115

# ALC row 123
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H123").Value = 52499.5
$ws.Range("J123").Value = 52499.5
$ws.Range("L123").Value = 52499.5
$ws.Range("N123").Value = -62299.5

# ALC row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 7414752
$ws.Range("I132").Value = 13895594
$ws.Range("K132").Value = 41686782
$ws.Range("M132").Value = -41684252

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6049.4717
$ws.Range("I32").Value = 5068.2554
$ws.Range("K32").Value = 5068.2554
$ws.Range("M32").Value = -4781.2554

# ARM row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1058
$ws.Range("I45").Value = 954.4211
$ws.Range("K45").Value = 954.4211
$ws.Range("M45").Value = -577.4211

# ARM row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1054.8077
$ws.Range("I122").Value = 1054.8077
$ws.Range("K122").Value = 3164.4231
$ws.Range("M122").Value = -714.4231

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2998.375
$ws.Range("I132").Value = 2456.1
$ws.Range("J132").Value = 3902.1667
$ws.Range("K132").Value = 7368.299999999999
$ws.Range("L132").Value = 11706.5001
$ws.Range("M132").Value = -4838.299999999999
$ws.Range("N132").Value = -16766.5001

# BSM row 37
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H37").Value = 5500
$ws.Range("I37").Value = 1000
$ws.Range("K37").Value = 1000
$ws.Range("M37").Value = -863

# BSM row 46
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()

# BSM row 94
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 50000840
$ws.Range("I94").Value = 62500300
$ws.Range("J94").Value = 3010
$ws.Range("K94").Value = 62500300
$ws.Range("L94").Value = 3010
$ws.Range("M94").Value = -62499849
$ws.Range("N94").Value = -3912

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 6097.8184
$ws.Range("I134").Value = 1038.4117
$ws.Range("K134").Value = 3115.2351
$ws.Range("M134").Value = -580.2351000000003

# CUL row 39
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 3005.4707
$ws.Range("J39").Value = 3146.2
$ws.Range("L39").Value = 9438.599999999999
$ws.Range("N39").Value = -10026.6

# CUL row 54
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()

# CUL row 55
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 2616.6667
$ws.Range("J55").Value = 3155.5557
$ws.Range("L55").Value = 9466.667099999999
$ws.Range("N55").Value = -9820.667099999999

# CUL row 74
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H74").Value = 5000
$ws.Range("I74").Value = 2000
$ws.Range("J74").Value = 5750
$ws.Range("K74").Value = 6000
$ws.Range("L74").Value = 17250
$ws.Range("M74").Value = -4939
$ws.Range("N74").Value = -19372

# CUL row 77
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H77").Value = 5000
$ws.Range("I77").Value = 2000
$ws.Range("J77").Value = 5750
$ws.Range("K77").Value = 18000
$ws.Range("L77").Value = 51750
$ws.Range("M77").Value = -12696
$ws.Range("N77").Value = -62358

# CUL row 81
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H81").Value = 2531
$ws.Range("I81").Value = 293
$ws.Range("J81").Value = 3650
$ws.Range("K81").Value = 879
$ws.Range("L81").Value = 10950
$ws.Range("M81").Value = 244
$ws.Range("N81").Value = -13196

# CUL row 84
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H84").Value = 2531
$ws.Range("I84").Value = 293
$ws.Range("J84").Value = 3650
$ws.Range("K84").Value = 2637
$ws.Range("L84").Value = 32850
$ws.Range("M84").Value = 2979
$ws.Range("N84").Value = -44082

# CUL row 87
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 2089.7144
$ws.Range("I87").Value = 882
$ws.Range("J87").Value = 3700
$ws.Range("K87").Value = 2646
$ws.Range("L87").Value = 11100
$ws.Range("M87").Value = -1398
$ws.Range("N87").Value = -13596

# CUL row 88
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H88").Value = 5926.316
$ws.Range("J88").Value = 6144.4443
$ws.Range("L88").Value = 18433.3329
$ws.Range("N88").Value = -19289.3329

# CUL row 90
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H90").Value = 2089.7144
$ws.Range("I90").Value = 882
$ws.Range("J90").Value = 3700
$ws.Range("K90").Value = 7938
$ws.Range("L90").Value = 33300
$ws.Range("M90").Value = -1698
$ws.Range("N90").Value = -45780

# CUL row 91
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H91").Value = 5926.316
$ws.Range("J91").Value = 6144.4443
$ws.Range("L91").Value = 18433.3329
$ws.Range("N91").Value = -21397.3329

# CUL row 98
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 332.33334
$ws.Range("I98").Value = 290.08334
$ws.Range("K98").Value = 870.2500200000001
$ws.Range("M98").Value = 627.7499799999999

# CUL row 113
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 683.8182
$ws.Range("I113").Value = 650
$ws.Range("J113").Value = 685.4286
$ws.Range("K113").Value = 1950
$ws.Range("L113").Value = 2056.2858
$ws.Range("M113").Value = 220
$ws.Range("N113").Value = -6396.2858

# CUL row 125
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H125").Value = 5005
$ws.Range("I125").Value = 1015
$ws.Range("J125").Value = 7000
$ws.Range("K125").Value = 3045
$ws.Range("L125").Value = 21000
$ws.Range("M125").Value = 1875
$ws.Range("N125").Value = -30840

# CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 22225968
$ws.Range("I131").Value = 125000310
$ws.Range("J131").Value = 4488
$ws.Range("K131").Value = 375000930
$ws.Range("L131").Value = 13464
$ws.Range("M131").Value = -374995890
$ws.Range("N131").Value = -23544

# CUL row 137
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 25003420
$ws.Range("J137").Value = 3870.2917
$ws.Range("L137").Value = 11610.8751
$ws.Range("N137").Value = -21810.8751

# GSM row 95
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()

# GSM row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1683.619
$ws.Range("I122").Value = 1515.0588
$ws.Range("J122").Value = 2400
$ws.Range("K122").Value = 4545.1764
$ws.Range("L122").Value = 7200
$ws.Range("M122").Value = -2095.1764
$ws.Range("N122").Value = -12100

# GSM row 126
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 1747.4117
$ws.Range("I126").Value = 1580.4
$ws.Range("K126").Value = 4741.200000000001
$ws.Range("M126").Value = -2271.200000000001

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 5891.3335
$ws.Range("I132").Value = 8277.5
$ws.Range("J132").Value = 3982.4
$ws.Range("K132").Value = 24832.5
$ws.Range("L132").Value = 11947.2
$ws.Range("M132").Value = -22302.5
$ws.Range("N132").Value = -17007.2

# LTW row 7
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2160.2
$ws.Range("I7").Value = 1767.1666
$ws.Range("K7").Value = 1767.1666
$ws.Range("M7").Value = -1655.1666

# LTW row 40
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2730.0386
$ws.Range("I40").Value = 1927.7142
$ws.Range("J40").Value = 6099.8
$ws.Range("K40").Value = 1927.7142
$ws.Range("L40").Value = 6099.8
$ws.Range("M40").Value = -1791.7142
$ws.Range("N40").Value = -6371.8

# LTW row 46
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 5335.5557
$ws.Range("I46").Value = 680
$ws.Range("K46").Value = 680
$ws.Range("M46").Value = -492

# LTW row 126
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 2160.2
$ws.Range("I126").Value = 1767.1666
$ws.Range("K126").Value = 5301.4998
$ws.Range("M126").Value = -2831.4998

# WVR row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 8066268
$ws.Range("I122").Value = 8930350
$ws.Range("K122").Value = 26791050
$ws.Range("M122").Value = -26788600

# WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1954.6538
$ws.Range("I132").Value = 1705.7142
$ws.Range("J132").Value = 3000.2
$ws.Range("K132").Value = 5117.142599999999
$ws.Range("L132").Value = 9000.599999999999
$ws.Range("M132").Value = -2587.142599999999
$ws.Range("N132").Value = -14060.6
